$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.106.06"
$ws.Range("E2").Value = "  -1.09%  "
$ws.Range("D3").Value = "1.787.02"
$ws.Range("E3").Value = "  -3.30%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.78"
$ws.Range("E5").Value = "  -0.96%  "
$ws.Range("E6").Value = "  -2.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.79"
$ws.Range("E8").Value = "  +1.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.285"
$ws.Range("E9").Value = "  -3.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0706"
$ws.Range("E10").Value = "  -2.64%  "
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").Value = "2.043.05"
$ws.Range("E12").Value = "  -3.31%  "
$ws.Range("D13").Value = "1.784.66"
$ws.Range("E13").Value = "  -3.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.82"
$ws.Range("E14").Value = "  -1.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.623"
$ws.Range("E15").Value = "  -4.56%  "
$ws.Range("D16").Value = "34.058.60"
$ws.Range("E16").Value = "  -1.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.16"
$ws.Range("E17").Value = "  -6.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.85"
$ws.Range("E18").Value = "  -3.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.18"
$ws.Range("E19").Value = "  -3.68%  "
$ws.Range("D20").Value = "0.0₃0785"
$ws.Range("E20").Value = "  -2.98%  "
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.80"
$ws.Range("E22").Value = "  -4.74%  "
$ws.Range("E23").Value = "  -5.14%  "
$ws.Range("E24").Value = "  -3.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.48"
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.33"
$ws.Range("E26").Value = "  -3.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.05"
$ws.Range("E27").Value = "  -3.53%  "
$ws.Range("E28").Value = "  -3.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0513"
$ws.Range("E30").Value = "  -5.00%  "
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("E32").Value = "  -4.88%  "
$ws.Range("E33").Value = "  -4.37%  "
$ws.Range("E34").Value = "  -7.02%  "
$ws.Range("D35").Value = "1.396.14"
$ws.Range("E35").Value = "  -4.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.642"
$ws.Range("E36").Value = "  -1.90%  "
$ws.Range("E37").Value = "  -2.23%  "
$ws.Range("E38").Value = "  -4.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.22"
$ws.Range("E39").Value = "  +2.77%  "
$ws.Range("E40").Value = "  -0.43%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.914"
$ws.Range("E41").Value = "  -6.60%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.70"
$ws.Range("E42").Value = "  -3.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "78.06"
$ws.Range("E43").Value = "  -5.94%  "
$ws.Range("D44").Value = "0.0₆0143"
$ws.Range("E44").Value = "  +12.44%  "
$ws.Range("E45").Value = "  +1.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.55"
$ws.Range("E46").Value = "  +1.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0498"
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "107.76"
$ws.Range("E48").Value = "  +0.65%  "
$ws.Range("E49").Value = "  -4.78%  "
$ws.Range("D50").Value = "1.944.06"
$ws.Range("E50").Value = "  -3.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.998"
$ws.Range("E51").Value = "  -0.40%  "
